# Refresh the cryptos list: updated prices / 1h volume %, and a ranking
# reshuffle among rows 47-49 (RenderToken overtakes VeChain, which in turn
# overtakes Maker). Values that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the source
# sheet, where the Price/Volume columns are plain strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.674.27"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "2.651.61"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'538.06"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'146.49"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").Value = "'6.75"
$ws.Range("E9").Value = "  +4.59%  "
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "3.116.42"
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "59.583.49"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "2.658.18"
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").Value = "'340.21"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").Value = "'4.42"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").Value = "'10.34"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'66.67"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").Value = "'0.417"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").Value = "0.0₃0749"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D31").Value = "'5.87"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").Value = "'18.93"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'150.94"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("D36").Value = "'0.840"
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("E37").Value = "  +0.74%  "
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  +1.73%  "
$ws.Range("D40").Value = "'287.26"
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "'0.608"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").Value = "'19.32"
$ws.Range("E45").Value = "  +4.10%  "
$ws.Range("D46").Value = "'0.0947"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'4.68"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0227"
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "1.968.73"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "'18.39"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").Value = "'111.10"
$ws.Range("E51").Value = "  +0.17%  "
